$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Last Install Date" marker ---
$ws.Range("O1").Value = 45449

# --- Row 6: CrossCode ---
$ws.Range("B6").Value = "0.5.0-pre2"
$ws.Range("C6").Value = 45457
$ws.Range("D6").Value = 45464

# --- Row 19: Shahrazad ---
$ws.Range("B19").Value = "0.1.1"
$ws.Range("C19").Value = 45457
$ws.Range("D19").Value = 45464

# --- Row 11: Jak and Daxter ---
$ws.Range("B11").Value = "0.0.4"
$ws.Range("C11").Value = 45453
$ws.Range("D11").Value = 45464

# --- Row 16: Outer Wilds ---
$ws.Range("B16").Value = "0.2.3"
$ws.Range("C16").Value = 45453
$ws.Range("D16").Value = 45464

# --- Row 8: Final Fantasy 5 Career Day (add Release Date) ---
$ws.Range("C9").Copy()
$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("C8").Value = 45417
$excel.CutCopyMode = 0

# --- Row 9: Final Fantasy 6 Worlds Collide ---
$ws.Range("C9").Value = 45453
$ws.Range("D9").Value = 45464

# --- Row 14: Minit (add Release Date) ---
$ws.Range("C9").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("C14").Value = 45405
$excel.CutCopyMode = 0

# --- Sheet view: update the active cell selection ---
$ws.Range("C10").Select()

# --- Workbook view window position/size (best effort; engine may not persist) ---
$win = $wb.Windows.Item(1)
$win.Left = 9660
$win.Top = 21480
$win.Width = 29040
$win.Height = 15840
